$p = $ppt.ActivePresentation
$s = $p.Slides.Item(23)
$shp = $s.Shapes.Item(1)

# 1) Nudge the code-listing rectangle up slightly (y offset 0 -> -76200)
$shp.Left = 228600
$shp.Top = -76200

$tr = $shp.TextFrame.TextRange

# 2) Move the "qs[i] = sum / length(range);" assignment out of the inner
#    "for (j in range)" loop: close that loop ("\t\t}") right after the
#    "sjs[j]" line, then de-indent the assignment by one tab and drop the
#    old trailing "\t\t}" paragraph that used to close the inner loop.

# Paragraph with "...sjs[j]" - append a tab and split the inner loop's
# closing brace into its own new paragraph right after it.
$sumLine = $tr.Paragraphs(29, 1)
$sumLine.InsertAfter("`r`t`t}")

# The assignment line is now paragraph 31; replace its leading three tabs
# with two, leaving the rest of the statement untouched.
$assignLine = $tr.Paragraphs(31, 1)
$assignLine.Text = "`t`tqs[i] = sum / length(range);`r"

# The paragraph that used to close the inner loop right after the
# assignment ("\t\t}") is now redundant (its brace already got placed
# above) - delete it.
$oldCloser = $tr.Paragraphs(32, 1)
$oldCloser.Text = ""
$oldCloser.Delete()
